$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.423.18'
$ws.Range("E2").Value = '  +2.17%  '

$ws.Range("D3").Value = '3.697.23'
$ws.Range("E3").Value = '  +1.59%  '

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.28%  '

$ws.Range("D5").Value = "'611.51"
$ws.Range("E5").Value = '  +6.63%  '

$ws.Range("D6").Value = "'194.82"
$ws.Range("E6").Value = '  +15.84%  '

$ws.Range("E7").Value = '  +3.31%  '

$ws.Range("E8").Value = '  -0.21%  '

$ws.Range("D9").Value = "'0.723"
$ws.Range("E9").Value = '  +4.58%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = "'0.160"
$ws.Range("E10").Value = '  +1.13%  '

$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").Value = "'59.57"
$ws.Range("E11").Value = '  +18.55%  '

$ws.Range("E12").Value = '  +1.37%  '

$ws.Range("D13").Value = "'10.42"
$ws.Range("E13").Value = '  +1.54%  '

$ws.Range("D14").Value = '4.295.29'
$ws.Range("E14").Value = '  +1.30%  '

$ws.Range("D15").Value = '3.706.92'
$ws.Range("E15").Value = '  +0.88%  '

$ws.Range("D16").Value = "'19.41"
$ws.Range("E16").Value = '  +1.88%  '

$ws.Range("E17").Value = '  +1.28%  '

$ws.Range("D18").Value = "'1.14"
$ws.Range("E18").Value = '  +4.28%  '

$ws.Range("D19").Value = "'12.81"
$ws.Range("E19").Value = '  +1.26%  '

$ws.Range("D20").Value = '68.324.64'
$ws.Range("E20").Value = '  +2.14%  '

$ws.Range("D21").Value = "'408.09"
$ws.Range("E21").Value = '  +2.12%  '

$ws.Range("E22").Value = '  +4.37%  '

$ws.Range("D23").Value = "'89.88"
$ws.Range("E23").Value = '  +3.99%  '

$ws.Range("D24").Value = "'11.50"
$ws.Range("E24").Value = '  +9.47%  '

$ws.Range("E25").Value = '  +2.78%  '

$ws.Range("D26").Value = "'13.04"
$ws.Range("E26").Value = '  +4.18%  '

$ws.Range("E27").Value = '  +1.14%  '

$ws.Range("D28").Value = "'3.76"
$ws.Range("E28").Value = '  +2.48%  '

$ws.Range("D29").Value = "'9.59"
$ws.Range("E29").Value = '  +3.42%  '

$ws.Range("D30").Value = "'32.67"
$ws.Range("E30").Value = '  +1.94%  '

$ws.Range("D31").Value = "'7.77"
$ws.Range("E31").Value = '  +5.55%  '

$ws.Range("D32").Value = "'48.00"
$ws.Range("E32").Value = '  +13.35%  '

$ws.Range("E33").Value = '  +3.93%  '

$ws.Range("D34").Value = "'0.121"
$ws.Range("E34").Value = '  +6.38%  '

$ws.Range("D35").Value = "'633.73"
$ws.Range("E35").Value = '  +8.97%  '

$ws.Range("D36").Value = "'67.48"
$ws.Range("E36").Value = '  +5.06%  '

$ws.Range("E37").Value = '  +6.52%  '

$ws.Range("D38").Value = '0.0₃0817'
$ws.Range("E38").Value = '  -6.26%  '

$ws.Range("E39").Value = '  +0.19%  '

$ws.Range("E40").Value = '  -0.21%  '

$ws.Range("E41").Value = '  +5.90%  '

$ws.Range("E42").Value = '  +3.69%  '

$ws.Range("E43").Value = '  +3.94%  '

$ws.Range("E44").Value = '  +0.91%  '

$ws.Range("D45").Value = '2.931.63'
$ws.Range("E45").Value = '  +6.11%  '

$ws.Range("E46").Value = '  +6.49%  '

$ws.Range("D47").Value = "'9.32"
$ws.Range("E47").Value = '  +3.44%  '

$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = "'146.07"
$ws.Range("E48").Value = '  +3.72%  '

$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").Value = "'2.69"
$ws.Range("E49").Value = '  +2.90%  '

$ws.Range("D50").Value = "'2.68"
$ws.Range("E50").Value = '  -4.51%  '

$ws.Range("D51").Value = "'3.04"
$ws.Range("E51").Value = '  -2.11%  '
